# Expand groupby functionality to include multi-field keys:
# Add a second worksheet "groupby_2" (positioned right after "groupby")
# containing an expanded version of the groupby dataset with an extra
# numeric column, and make it the active sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after the existing "groupby" sheet.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "groupby_2"

# Header-ish value (same pattern as sheet1, shifted one column right).
$ws2.Range("D1").Value = 2018

$ws2.Range("A2").Value = "a"
$ws2.Range("B2").Value = "b"
$ws2.Range("C2").Value = 1
$ws2.Range("D2").Value = 1

$ws2.Range("A3").Value = "a"
$ws2.Range("B3").Value = "c"
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 2

$ws2.Range("A4").Value = "a"
$ws2.Range("B4").Value = "d"
$ws2.Range("C4").Value = 2
$ws2.Range("D4").Value = 3

$ws2.Range("A5").Value = "b"
$ws2.Range("B5").Value = "b"
$ws2.Range("C5").Value = 3
$ws2.Range("D5").Value = 4

$ws2.Range("A6").Value = "c"
$ws2.Range("B6").Value = "b"
$ws2.Range("C6").Value = 3
$ws2.Range("D6").Value = 5

$ws2.Range("A7").Value = "a"
$ws2.Range("B7").Value = "b"
$ws2.Range("C7").Value = 2
$ws2.Range("D7").Value = 6

$ws2.Range("A8").Value = "a"
$ws2.Range("B8").Value = "c"
$ws2.Range("C8").Value = 2
$ws2.Range("D8").Value = 7

$ws2.Range("A9").Value = "a"
$ws2.Range("B9").Value = "d"
$ws2.Range("C9").Value = 3
$ws2.Range("D9").Value = 8

$ws2.Range("A10").Value = "a"
$ws2.Range("B10").Value = "c"
$ws2.Range("C10").Value = 3
$ws2.Range("D10").Value = 9

# Reset the selection on the original sheet back to A1 and make sure it
# is no longer the tab-selected sheet.
[void]$ws1.Range("A1").Select()

# Select A6 on the new sheet and leave it as the active sheet/tab.
[void]$ws2.Range("A6").Select()
